$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6534.3438
$ws.Range("I43").Value = 5323.353
$ws.Range("J43").Value = 7906.8
$ws.Range("K43").Value = 5323.353
$ws.Range("L43").Value = 7906.8
$ws.Range("M43").Value = -5254.353
$ws.Range("N43").Value = -8044.8
$ws.Range("H76").Value = 6675625
$ws.Range("I76").Value = 9100198
$ws.Range("J76").Value = 8048.5
$ws.Range("K76").Value = 9100198
$ws.Range("L76").Value = 8048.5
$ws.Range("M76").Value = -9099883
$ws.Range("N76").Value = -8678.5
$ws.Range("H79").Value = 6675625
$ws.Range("I79").Value = 9100198
$ws.Range("J79").Value = 8048.5
$ws.Range("K79").Value = 9100198
$ws.Range("L79").Value = 8048.5
$ws.Range("M79").Value = -9099106
$ws.Range("N79").Value = -10232.5
$ws.Range("H86").Value = 2886.9092
$ws.Range("I86").Value = 1965.2858
$ws.Range("K86").Value = 1965.2858
$ws.Range("M86").Value = -842.2858000000001
$ws.Range("H89").Value = 2886.9092
$ws.Range("I89").Value = 1965.2858
$ws.Range("K89").Value = 9826.429
$ws.Range("M89").Value = -4210.429
$ws.Range("H103").Value = 809.9091
$ws.Range("I103").Value = 821.8570999999999
$ws.Range("J103").Value = 789
$ws.Range("K103").Value = 2465.5713
$ws.Range("L103").Value = 2367
$ws.Range("M103").Value = -1879.5713
$ws.Range("N103").Value = -3539
$ws.Range("H138").Value = 3956.875
$ws.Range("I138").Value = 3368.2666
$ws.Range("J138").Value = 4172.2197
$ws.Range("K138").Value = 10104.7998
$ws.Range("L138").Value = 12516.6591
$ws.Range("M138").Value = -4964.799800000001
$ws.Range("N138").Value = -22796.6591

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1737146.8
$ws.Range("I2").Value = 2045790.9
$ws.Range("K2").Value = 2045790.9
$ws.Range("M2").Value = -2045677.9
$ws.Range("H32").Value = 20285.674
$ws.Range("I32").Value = 19749.377
$ws.Range("K32").Value = 19749.377
$ws.Range("M32").Value = -19462.377
$ws.Range("H61").Value = 8134993.5
$ws.Range("I61").Value = 11115046
$ws.Range("J61").Value = 7578.091
$ws.Range("K61").Value = 11115046
$ws.Range("L61").Value = 7578.091
$ws.Range("M61").Value = -11114834
$ws.Range("N61").Value = -8002.091
$ws.Range("H116").Value = 1737146.8
$ws.Range("I116").Value = 2045790.9
$ws.Range("K116").Value = 2045790.9
$ws.Range("M116").Value = -2043496.9
$ws.Range("H136").Value = 8134993.5
$ws.Range("I136").Value = 11115046
$ws.Range("J136").Value = 7578.091
$ws.Range("K136").Value = 33345138
$ws.Range("L136").Value = 22734.273
$ws.Range("M136").Value = -33342588
$ws.Range("N136").Value = -27834.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1737146.8
$ws.Range("I3").Value = 2045790.9
$ws.Range("K3").Value = 2045790.9
$ws.Range("M3").Value = -2045676.9
$ws.Range("H20").Value = 3628.0312
$ws.Range("I20").Value = 3325
$ws.Range("K20").Value = 3325
$ws.Range("M20").Value = -3078
$ws.Range("H96").Value = 5808.6665
$ws.Range("I96").Value = 5808.6665
$ws.Range("K96").Value = 5808.6665
$ws.Range("M96").Value = -3062.6665
$ws.Range("H105").Value = 45467930
$ws.Range("I105").Value = 66684690
$ws.Range("J105").Value = 3433.1428
$ws.Range("K105").Value = 66684690
$ws.Range("L105").Value = 3433.1428
$ws.Range("M105").Value = -66682943
$ws.Range("N105").Value = -6927.1428
$ws.Range("H117").Value = 112999.5
$ws.Range("J117").Value = 112999.5
$ws.Range("L117").Value = 112999.5
$ws.Range("N117").Value = -122177.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 46296.668
$ws.Range("J20").Value = 46296.668
$ws.Range("L20").Value = 46296.668
$ws.Range("N20").Value = -46768.668
$ws.Range("H30").Value = 46296.668
$ws.Range("J30").Value = 46296.668
$ws.Range("L30").Value = 46296.668
$ws.Range("N30").Value = -46478.668
$ws.Range("H31").Value = 32263250
$ws.Range("J31").Value = 8720
$ws.Range("L31").Value = 8720
$ws.Range("N31").Value = -9310
$ws.Range("H34").Value = 32263250
$ws.Range("J34").Value = 8720
$ws.Range("L34").Value = 8720
$ws.Range("N34").Value = -9124
$ws.Range("H105").Value = 2997
$ws.Range("I105").Value = 2997
$ws.Range("K105").Value = 2997
$ws.Range("M105").Value = -1250
$ws.Range("H106").Value = 54499.75
$ws.Range("J106").Value = 54499.75
$ws.Range("L106").Value = 54499.75
$ws.Range("N106").Value = -57023.75
$ws.Range("H107").Value = 62500744
$ws.Range("I107").Value = 76923610
$ws.Range("K107").Value = 76923610
$ws.Range("M107").Value = -76921690
$ws.Range("H128").Value = 46296.668
$ws.Range("J128").Value = 46296.668
$ws.Range("L128").Value = 46296.668
$ws.Range("N128").Value = -56256.668
$ws.Range("H141").Value = 208989.25
$ws.Range("J141").Value = 225543.61
$ws.Range("L141").Value = 225543.61
$ws.Range("N141").Value = -235903.61

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 799.0833
$ws.Range("I14").Value = 799.0833
$ws.Range("K14").Value = 2397.2499
$ws.Range("M14").Value = -2224.2499
$ws.Range("H61").Value = 111.28125
$ws.Range("I61").Value = 49.296295
$ws.Range("J61").Value = 446
$ws.Range("K61").Value = 147.888885
$ws.Range("L61").Value = 1338
$ws.Range("M61").Value = 67.11111499999998
$ws.Range("N61").Value = -1768
$ws.Range("H62").Value = 1396.5385
$ws.Range("I62").Value = 1242.2972
$ws.Range("K62").Value = 3726.8916
$ws.Range("M62").Value = -3040.8916
$ws.Range("H63").Value = 4537.8
$ws.Range("I63").Value = 4537.8
$ws.Range("K63").Value = 13613.4
$ws.Range("M63").Value = -12864.4
$ws.Range("H64").Value = 4651.5
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15540
$ws.Range("H65").Value = 1396.5385
$ws.Range("I65").Value = 1242.2972
$ws.Range("K65").Value = 11180.6748
$ws.Range("M65").Value = -7748.674800000001
$ws.Range("H66").Value = 4537.8
$ws.Range("I66").Value = 4537.8
$ws.Range("K66").Value = 40840.2
$ws.Range("M66").Value = -37096.2
$ws.Range("H67").Value = 4651.5
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16872
$ws.Range("H107").Value = 1815.8182
$ws.Range("J107").Value = 1815.8182
$ws.Range("L107").Value = 5447.4546
$ws.Range("N107").Value = -9287.454600000001
$ws.Range("H137").Value = 2737.1052
$ws.Range("J137").Value = 5714.143
$ws.Range("L137").Value = 17142.429
$ws.Range("N137").Value = -27342.429
$ws.Range("H138").Value = 403562.16
$ws.Range("I138").Value = 2404.2222
$ws.Range("J138").Value = 629213.5
$ws.Range("K138").Value = 7212.6666
$ws.Range("L138").Value = 1887640.5
$ws.Range("M138").Value = -2072.6666
$ws.Range("N138").Value = -1897920.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28901.285
$ws.Range("I2").Value = 368.25
$ws.Range("K2").Value = 368.25
$ws.Range("M2").Value = -255.25
$ws.Range("H12").Value = 772
$ws.Range("I12").Value = 176.5
$ws.Range("J12").Value = 1367.5
$ws.Range("K12").Value = 176.5
$ws.Range("L12").Value = 1367.5
$ws.Range("M12").Value = -36.5
$ws.Range("N12").Value = -1647.5
$ws.Range("H113").Value = 2563
$ws.Range("I113").Value = 1745.9166
$ws.Range("K113").Value = 1745.9166
$ws.Range("M113").Value = 424.0834
$ws.Range("H122").Value = 3620.3333
$ws.Range("J122").Value = 5134.3687
$ws.Range("L122").Value = 15403.1061
$ws.Range("N122").Value = -20303.1061
$ws.Range("H132").Value = 3886.75
$ws.Range("I132").Value = 3054.5312
$ws.Range("K132").Value = 9163.5936
$ws.Range("M132").Value = -6633.5936

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2242.2188
$ws.Range("I61").Value = 2242.2188
$ws.Range("K61").Value = 2242.2188
$ws.Range("M61").Value = -2040.2188
$ws.Range("H100").Value = 9261397
$ws.Range("I100").Value = 22728348
$ws.Range("K100").Value = 22728348
$ws.Range("M100").Value = -22727807
$ws.Range("H113").Value = 2242.2188
$ws.Range("I113").Value = 2242.2188
$ws.Range("K113").Value = 2242.2188
$ws.Range("M113").Value = -72.2188000000001
$ws.Range("H127").Value = 150000
$ws.Range("J127").Value = 150000
$ws.Range("L127").Value = 150000
$ws.Range("N127").Value = -159920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8829.326999999999
$ws.Range("I81").Value = 3710.2593
$ws.Range("K81").Value = 7420.5186
$ws.Range("M81").Value = -6359.5186
$ws.Range("H84").Value = 8829.326999999999
$ws.Range("I84").Value = 3710.2593
$ws.Range("K84").Value = 37102.593
$ws.Range("M84").Value = -31798.593
$ws.Range("H110").Value = 72999
$ws.Range("J110").Value = 72999
$ws.Range("L110").Value = 72999
$ws.Range("N110").Value = -81179
$ws.Range("H116").Value = 116398.75
$ws.Range("J116").Value = 116398.75
$ws.Range("L116").Value = 116398.75
$ws.Range("N116").Value = -125576.75
